$p = $ppt.ActivePresentation

$oldDate = "15/02/2025"
$newDate = "17/02/2025"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*" -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master footer date placeholder.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every custom (slide) layout's footer date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($l = 1; $l -le $layouts.Count; $l++) {
    $layout = $layouts.Item($l)
    Update-DatePlaceholder $layout.Shapes
}
